# chl-transect-info.xlsx revision:
#  - revised abstract, methods (metadata-only changes not reflected in cell data)
#  - added ORCID for Kate Morkeski (Personnel sheet, userId column F, row 13)
#  - additional_info -> active sheet/selection moved to Personnel sheet at G17

$wb = $excel.ActiveWorkbook

# "Personnel" is the 3rd sheet (rId3 -> sheet3.xml)
$ws = $wb.Worksheets.Item(3)

# Add Kate Morkeski's ORCID iD to the userId column (F13)
$ws.Range("F13").Value = "0000-0002-2903-5851"

# Make Personnel the active sheet, with G17 selected (mirrors the saved view state)
$ws.Activate()
$ws.Range("G17").Select() | Out-Null
